$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. LoginDetails: reorder hyperlinks (A3, B3, B2, A2) and fix the selection.
# ---------------------------------------------------------------------------
$login = $wb.Worksheets.Item("LoginDetails")

$login.Hyperlinks.Delete()
$login.Hyperlinks.Add($login.Range("A3"), "mailto:tushar.jadhav.work@gmail.com")
$login.Hyperlinks.Add($login.Range("B3"), "mailto:Jadhav@1228")
$login.Hyperlinks.Add($login.Range("B2"), "mailto:Sakshi@1228")
$login.Hyperlinks.Add($login.Range("A2"), "mailto:tusharjadhav228@gmail.com")
$login.Range("A2:B3").Style = "Hyperlink"

# ---------------------------------------------------------------------------
# 2. Insert a new "UserDetails" sheet between LoginDetails and
#    PassengerDetails, populated with a copy of PassengerDetails' row 2
#    (only that single row - no header, no extra rows/cols).
# ---------------------------------------------------------------------------
$user = $wb.Worksheets.Add($wb.Worksheets.Item("PassengerDetails"))
$user.Name = "UserDetails"

# The Worksheets collection shifts once the new sheet is inserted, so the
# PassengerDetails handle must be re-acquired after the Add() call.
$passenger = $wb.Worksheets.Item("PassengerDetails")

$user.Range("A2").Value = "tusharjadhav228@gmail.com"
$user.Range("B2").Value = "Sakshi@1228"
$user.Range("C2").Value = "Tushar"
$user.Range("D2").Value = "Jadhav"
$user.Range("E2").Value = 3
$user.Range("F2").Value = 8
$user.Range("G2").Value = 1999
$user.Range("H2").Value = "India"
$user.Range("I2").Value = "Asalfa"
$user.Range("J2").Value = "Mumbai"
$user.Range("K2").Value = 8745963287
$user.Range("L2").Value = "tusharjadhav123@gmail.com"
$user.Range("M2").Value = "4111 1111 1111 1111"
$user.Range("N2").Value = "Tushar Jadhav"
$user.Range("O2").NumberFormat = "@"
$user.Range("O2").Value = "08"
$user.Range("P2").Value = 2026
$user.Range("Q2").Value = 123

$user.Hyperlinks.Add($user.Range("A2"), "mailto:tusharjadhav228@gmail.com")
$user.Hyperlinks.Add($user.Range("B2"), "mailto:Sakshi@1228")
$user.Range("A2:B2").Style = "Hyperlink"

# ---------------------------------------------------------------------------
# 3. PassengerDetails: update the selection (whole row 2 selected).
# ---------------------------------------------------------------------------
$passenger.Range("A2:XFD2").Select()

# ---------------------------------------------------------------------------
# 4. Selections + active sheet/tab.
#    LoginDetails gets cell C22 selected, UserDetails gets C6 selected and
#    becomes the active (visible) tab, matching activeTab="1".
# ---------------------------------------------------------------------------
$login.Range("C22").Select()
$user.Range("C6").Select()
$user.Activate()
